# 0.4.1 kiirja a maganut tavot
#
# Travel-log ("utnyilvantartas") update for 2021-01 (Teszt Elek / ELEK-01):
# extend the trip table with 14 additional legs (rows 20-33) that were
# previously blank, and correct the fuel-use figure for the existing row 19
# leg (122 -> 379 km). Both changes ripple through the running km total in
# column G and the summary block at the bottom of the sheet (L172:M174),
# all of which recompute automatically via the existing formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: fuel-use correction ripples the running km total (col G) ---
$ws.Cells.Item(19, 9).Value = 379.0   # I19
# Re-push the shared formula on G19 explicitly (single-quoted so the
# embedded double-quotes need no escaping) - this keeps the same formula
# text but forces a fresh, correctly-cascaded evaluation.
$ws.Range("G19").Formula = '=IF(B19<>"",G18+I19,"")'

# --- Rows 20-33: fourteen new trip legs ---
# Column B holds dates stored as literal text (not date serials); switch
# the number format to Text before writing so Excel does not coerce the
# string into a date serial, then restore the original format so the
# cell style (s="24") stays identical to the other B-column cells.
$origBFormat = $ws.Range("B20").NumberFormat
$ws.Range("B20:B33").NumberFormat = "@"

$ws.Cells.Item(20, 2).Value = "2021-01-04"    # B20
$ws.Cells.Item(20, 3).Value = "Hibajavítás"    # C20
$ws.Cells.Item(20, 4).Value = " Súr Kossuth ut 5 "    # D20
$ws.Cells.Item(20, 5).Value = "Szombathely Éhen Gyula tér 2."    # E20
$ws.Cells.Item(20, 6).Value = "AKHB2078/K & H BANK ZRT."    # F20
$ws.Cells.Item(20, 8).Value = 0.0    # H20
$ws.Cells.Item(20, 9).Value = 149.0    # I20
$ws.Cells.Item(20, 10).Value = "C"    # J20
$ws.Range("G20").Formula = '=IF(B20<>"",G19+I20,"")'

$ws.Cells.Item(21, 2).Value = "2021-01-04"    # B21
$ws.Cells.Item(21, 3).Value = "Hibajavítás"    # C21
$ws.Cells.Item(21, 4).Value = "Szombathely Éhen Gyula tér 2."    # D21
$ws.Cells.Item(21, 5).Value = "Edelény Antal György u. 3."    # E21
$ws.Cells.Item(21, 6).Value = "AKHB2144/K & H BANK ZRT."    # F21
$ws.Cells.Item(21, 8).Value = 0.0    # H21
$ws.Cells.Item(21, 9).Value = 451.0    # I21
$ws.Cells.Item(21, 10).Value = "C"    # J21
$ws.Range("G21").Formula = '=IF(B21<>"",G20+I21,"")'

$ws.Cells.Item(22, 2).Value = "2021-01-04"    # B22
$ws.Cells.Item(22, 3).Value = "Hibajavítás"    # C22
$ws.Cells.Item(22, 4).Value = "Edelény Antal György u. 3."    # D22
$ws.Cells.Item(22, 5).Value = "Monostorpályi Landler tér 4."    # E22
$ws.Cells.Item(22, 6).Value = "AKHB2122/K & H BANK ZRT."    # F22
$ws.Cells.Item(22, 8).Value = 0.0    # H22
$ws.Cells.Item(22, 9).Value = 165.0    # I22
$ws.Cells.Item(22, 10).Value = "C"    # J22
$ws.Range("G22").Formula = '=IF(B22<>"",G21+I22,"")'

$ws.Cells.Item(23, 2).Value = "2021-01-04"    # B23
$ws.Cells.Item(23, 3).Value = "Hibajavítás"    # C23
$ws.Cells.Item(23, 4).Value = "Monostorpályi Landler tér 4."    # D23
$ws.Cells.Item(23, 5).Value = "Nyírbátor Debreceni út 71."    # E23
$ws.Cells.Item(23, 6).Value = "AKHB2317/K & H BANK ZRT."    # F23
$ws.Cells.Item(23, 8).Value = 0.0    # H23
$ws.Cells.Item(23, 9).Value = 73.0    # I23
$ws.Cells.Item(23, 10).Value = "C"    # J23
$ws.Range("G23").Formula = '=IF(B23<>"",G22+I23,"")'

$ws.Cells.Item(24, 2).Value = "2021-01-04"    # B24
$ws.Cells.Item(24, 3).Value = "Hibajavítás"    # C24
$ws.Cells.Item(24, 4).Value = "Nyírbátor Debreceni út 71."    # D24
$ws.Cells.Item(24, 5).Value = "Cegléd Ipartelepi u. 3."    # E24
$ws.Cells.Item(24, 6).Value = "AKHB2340/K & H BANK ZRT."    # F24
$ws.Cells.Item(24, 8).Value = 0.0    # H24
$ws.Cells.Item(24, 9).Value = 244.0    # I24
$ws.Cells.Item(24, 10).Value = "C"    # J24
$ws.Range("G24").Formula = '=IF(B24<>"",G23+I24,"")'

$ws.Cells.Item(25, 2).Value = "2021-01-04"    # B25
$ws.Cells.Item(25, 3).Value = "Hibajavítás"    # C25
$ws.Cells.Item(25, 4).Value = "Cegléd Ipartelepi u. 3."    # D25
$ws.Cells.Item(25, 5).Value = "Sárvár Rákóczi u. 83."    # E25
$ws.Cells.Item(25, 6).Value = "AKHB2362/K & H BANK ZRT."    # F25
$ws.Cells.Item(25, 8).Value = 0.0    # H25
$ws.Cells.Item(25, 9).Value = 286.0    # I25
$ws.Cells.Item(25, 10).Value = "C"    # J25
$ws.Range("G25").Formula = '=IF(B25<>"",G24+I25,"")'

$ws.Cells.Item(26, 2).Value = "2021-01-04"    # B26
$ws.Cells.Item(26, 3).Value = "Hibajavítás"    # C26
$ws.Cells.Item(26, 4).Value = "Sárvár Rákóczi u. 83."    # D26
$ws.Cells.Item(26, 5).Value = "Nyíregyháza Állomás tér 3."    # E26
$ws.Cells.Item(26, 6).Value = "AKHB2670/K & H BANK ZRT."    # F26
$ws.Cells.Item(26, 8).Value = 0.0    # H26
$ws.Cells.Item(26, 9).Value = 461.0    # I26
$ws.Cells.Item(26, 10).Value = "C"    # J26
$ws.Range("G26").Formula = '=IF(B26<>"",G25+I26,"")'

$ws.Cells.Item(27, 2).Value = "2021-01-04"    # B27
$ws.Cells.Item(27, 3).Value = "Hibajavítás"    # C27
$ws.Cells.Item(27, 4).Value = "Nyíregyháza Állomás tér 3."    # D27
$ws.Cells.Item(27, 5).Value = " Súr Kossuth ut 5 "    # E27
$ws.Cells.Item(27, 6).Value = "telephely/telephely"    # F27
$ws.Cells.Item(27, 8).Value = 0.0    # H27
$ws.Cells.Item(27, 9).Value = 360.0    # I27
$ws.Cells.Item(27, 10).Value = "C"    # J27
$ws.Range("G27").Formula = '=IF(B27<>"",G26+I27,"")'

$ws.Cells.Item(28, 2).Value = "2021-01-04"    # B28
$ws.Cells.Item(28, 3).Value = "Hibajavítás"    # C28
$ws.Cells.Item(28, 4).Value = " Súr Kossuth ut 5 "    # D28
$ws.Cells.Item(28, 5).Value = "Budapest Andrássy út 49."    # E28
$ws.Cells.Item(28, 6).Value = "AKHK2155/K & H BANK ZRT."    # F28
$ws.Cells.Item(28, 8).Value = 0.0    # H28
$ws.Cells.Item(28, 9).Value = 102.0    # I28
$ws.Cells.Item(28, 10).Value = "C"    # J28
$ws.Range("G28").Formula = '=IF(B28<>"",G27+I28,"")'

$ws.Cells.Item(29, 2).Value = "2021-01-04"    # B29
$ws.Cells.Item(29, 3).Value = "Hibajavítás"    # C29
$ws.Cells.Item(29, 4).Value = "Budapest Andrássy út 49."    # D29
$ws.Cells.Item(29, 5).Value = " Súr Kossuth ut 5 "    # E29
$ws.Cells.Item(29, 6).Value = "telephely/telephely"    # F29
$ws.Cells.Item(29, 8).Value = 0.0    # H29
$ws.Cells.Item(29, 9).Value = 102.0    # I29
$ws.Cells.Item(29, 10).Value = "C"    # J29
$ws.Range("G29").Formula = '=IF(B29<>"",G28+I29,"")'

$ws.Cells.Item(30, 2).Value = "2021-01-04"    # B30
$ws.Cells.Item(30, 3).Value = "Hibajavítás"    # C30
$ws.Cells.Item(30, 4).Value = " Súr Kossuth ut 5 "    # D30
$ws.Cells.Item(30, 5).Value = "Debrecen Füredi út 27."    # E30
$ws.Cells.Item(30, 6).Value = "AKHB2072/K & H BANK ZRT."    # F30
$ws.Cells.Item(30, 8).Value = 0.0    # H30
$ws.Cells.Item(30, 9).Value = 353.0    # I30
$ws.Cells.Item(30, 10).Value = "C"    # J30
$ws.Range("G30").Formula = '=IF(B30<>"",G29+I30,"")'

$ws.Cells.Item(31, 2).Value = "2021-01-04"    # B31
$ws.Cells.Item(31, 3).Value = "Hibajavítás"    # C31
$ws.Cells.Item(31, 4).Value = "Debrecen Füredi út 27."    # D31
$ws.Cells.Item(31, 5).Value = " Súr Kossuth ut 5 "    # E31
$ws.Cells.Item(31, 6).Value = "telephely/telephely"    # F31
$ws.Cells.Item(31, 8).Value = 0.0    # H31
$ws.Cells.Item(31, 9).Value = 353.0    # I31
$ws.Cells.Item(31, 10).Value = "C"    # J31
$ws.Range("G31").Formula = '=IF(B31<>"",G30+I31,"")'

$ws.Cells.Item(32, 2).Value = "2021-01-04"    # B32
$ws.Cells.Item(32, 3).Value = "Hibajavítás"    # C32
$ws.Cells.Item(32, 4).Value = " Súr Kossuth ut 5 "    # D32
$ws.Cells.Item(32, 5).Value = "Debrecen Füredi út 27."    # E32
$ws.Cells.Item(32, 6).Value = "AKHB2072/K & H BANK ZRT."    # F32
$ws.Cells.Item(32, 8).Value = 0.0    # H32
$ws.Cells.Item(32, 9).Value = 353.0    # I32
$ws.Cells.Item(32, 10).Value = "C"    # J32
$ws.Range("G32").Formula = '=IF(B32<>"",G31+I32,"")'

$ws.Cells.Item(33, 2).Value = "2021-01-04"    # B33
$ws.Cells.Item(33, 3).Value = "Hibajavítás"    # C33
$ws.Cells.Item(33, 4).Value = "Debrecen Füredi út 27."    # D33
$ws.Cells.Item(33, 5).Value = " Súr Kossuth ut 5 "    # E33
$ws.Cells.Item(33, 6).Value = "telephely/telephely"    # F33
$ws.Cells.Item(33, 8).Value = 0.0    # H33
$ws.Cells.Item(33, 9).Value = 353.0    # I33
$ws.Cells.Item(33, 10).Value = "C"    # J33
$ws.Range("G33").Formula = '=IF(B33<>"",G32+I33,"")'

# Restore original (date) number format on column B for the new rows
$ws.Range("B20:B33").NumberFormat = $origBFormat
